# The post about "「雲」سَحَابَة" (row 16) was removed from the sheet.
# Deleting the entire row shifts every row below it up by one, which
# reproduces the renumbering seen across rows 17-196 -> 16-195 in the
# diff, and updates the sheet dimension from A1:C196 to A1:C195.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(16).Delete()
